$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.722.07"
$ws.Range("E2").Value = "  +11.70%  "

$ws.Range("D3").Value = "1.683.14"
$ws.Range("E3").Value = "  +6.33%  "

$ws.Range("E4").Value = "  -0.45%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.10"
$ws.Range("E5").Value = "  +3.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9954"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3682"
$ws.Range("E7").Value = "  +2.10%  "

$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.94"
$ws.Range("E8").Value = "  +19.10%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3419"
$ws.Range("E9").Value = "  +2.94%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.160"
$ws.Range("E10").Value = "  +4.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07225"
$ws.Range("E11").Value = "  +4.72%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9969"
$ws.Range("E12").Value = "  -0.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.083"
$ws.Range("E13").Value = "  +5.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.11"
$ws.Range("E14").Value = "  +4.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.684"
$ws.Range("E15").Value = "  +3.00%  "

$ws.Range("D16").Value = "1.679.40"
$ws.Range("E16").Value = "  +5.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001101"
$ws.Range("E17").Value = "  +4.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9949"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06647"
$ws.Range("E19").Value = "  +1.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "80.70"
$ws.Range("E20").Value = "  +6.82%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.36"
$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.065"
$ws.Range("E22").Value = "  +3.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.07"
$ws.Range("E23").Value = "  +5.08%  "

$ws.Range("D24").Value = "24.646.52"
$ws.Range("E24").Value = "  +11.26%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.404"
$ws.Range("E25").Value = "  +1.88%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.654"
$ws.Range("E26").Value = "  +7.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.56"
$ws.Range("E27").Value = "  +2.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.42"
$ws.Range("E28").Value = "  +2.12%  "

$ws.Range("D29").Value = "1.866.41"
$ws.Range("E29").Value = "  +6.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.45"
$ws.Range("E30").Value = "  +5.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.241"
$ws.Range("E31").Value = "  +8.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.020"
$ws.Range("E32").Value = "  +2.92%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9739"
$ws.Range("E33").Value = "  +7.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08402"
$ws.Range("E34").Value = "  +3.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.683"
$ws.Range("E35").Value = "  +3.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.29"
$ws.Range("E36").Value = "  +6.07%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.296"
$ws.Range("E37").Value = "  +4.92%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06333"
$ws.Range("E38").Value = "  +6.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02302"
$ws.Range("E39").Value = "  +6.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.608"
$ws.Range("E40").Value = "  +4.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.243"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2081"
$ws.Range("E42").Value = "  +6.15%  "

$ws.Range("E43").Value = "  +6.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9950"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("E45").Value = "  +0.18%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.00"
$ws.Range("E46").Value = "  +2.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5870"
$ws.Range("E47").Value = "  +6.26%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.35"
$ws.Range("E48").Value = "  +1.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.996"
$ws.Range("E49").Value = "  +4.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07221"
$ws.Range("E50").Value = "  +7.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.59"
$ws.Range("E51").Value = "  +5.18%  "
